# Merge the trailing "console_draw)" run with the following "." run into a
# single run "console_draw)." (last paragraph of the document, just before
# the closing sectPr).
$d = $word.ActiveDocument

# Locate the last paragraph and find "console_draw)" within it.
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraStart = $p.Range.Start
$paraText = $p.Range.Text
$offset = $paraText.IndexOf("console_draw)")

$markerEnd = $paraStart + $offset + "console_draw)".Length
$dotEnd = $markerEnd + 1

# Range covering exactly the "." run that follows "console_draw)".
$dotRange = $d.Range($markerEnd, $dotEnd)

# Use Find/Replace (rather than a plain Range.Text assignment, which is a
# no-op for identical text) scoped tightly to that single-character range.
# This causes the runtime to re-coalesce it with the immediately preceding
# run (which ends in "console_draw)") -- picking up that run's formatting --
# without touching any earlier, differently-bounded runs in the paragraph.
$dotRange.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, ".", 2)
